# "finished tasks for importing monthly events and content pass"
#
# The OneTimeEvents sheet is reorganized into a ContentPass sheet: the
# name_en/name_zh columns (old B/C) move to the end of the row (after
# icon_url), so start_time/end_time (old D/E) become the new B/C and
# icon_url (old F) becomes the new D.  Final layout:
#   A name | B start_time | C end_time | D icon_url | E name_en | F name_zh

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OneTimeEvents")

# Move name_en (column B) to the end of the row (after icon_url, col F).
$ws.Columns.Item(2).Cut() | Out-Null
$ws.Columns.Item(7).Insert() | Out-Null

# Move name_zh (now column B again, since B shifted left) to the end too.
$ws.Columns.Item(2).Cut() | Out-Null
$ws.Columns.Item(7).Insert() | Out-Null

# icon_url (now column D) needs a wider column to fit its longer values.
$ws.Columns.Item(4).ColumnWidth = 37.142857142857146

# Reflect the author's final cursor position on the sheet.
$ws.Range("C8").Select() | Out-Null

# The sheet now represents the "content pass" data set.
$ws.Name = "ContentPass"
